$d = $word.ActiveDocument

# -------------------------------------------------------------------
# Edit 1: remove ", and report to jail on October 14, 2022, at 7:00
# p.m" from the sentencing paragraph, leaving "...by October 11,
# 2022." The trailing period must stay in its own (non-bold) run, so
# we copy it, delete the whole span (including the period, which is
# required for the runtime to also drop the orphaned spellStart/
# spellEnd proofErr markers around "p.m"), then paste the period back
# using the clipboard so its original run formatting is preserved
# exactly.
# -------------------------------------------------------------------

$afterDate = $d.Content
$afterDate.Find.Execute("October 11, 2022")
$startPos = $afterDate.End

$pmSearch = $d.Content
$pmSearch.Start = $startPos
$pmSearch.Find.Execute("p.m.")
$endPos = $pmSearch.End

$periodRange = $d.Range($endPos - 1, $endPos)
$periodRange.Copy()

$toDelete = $d.Range($startPos, $endPos)
$toDelete.Delete()

$insertionPoint = $d.Range($startPos, $startPos)
$insertionPoint.Paste()

# -------------------------------------------------------------------
# Edit 2: remove the whole "Restitution. The Defendant must pay
# restitution in the amount of $5,000 to Justin Kudela no later than
# October 11, 2022 in order to successfully complete the diversion
# program." sentence (plus its trailing line breaks), leaving just
# the pre-existing empty run immediately before "Fines and Costs."
# -------------------------------------------------------------------

$restStart = $d.Content
$restStart.Find.Execute("Restitution.")
$restitutionStart = $restStart.Start

$finesSearch = $d.Content
$finesSearch.Start = $restitutionStart
$finesSearch.Find.Execute("Fines and Costs.")
$finesStart = $finesSearch.Start

$restRange = $d.Range($restitutionStart, $finesStart)
$restRange.Delete()
